# Leave-card update (periodic refresh of the leave card table):
#  - Shift the PERIOD date (column A) for rows 170-212 from the 1st of
#    each month to the last day of that same month.
#  - Fill in newly-earned leave credits (1.25) for rows 175-179 in the
#    "EARNED" column (C) of the first half of Table1; the mirrored
#    "EARNED " column (G) is a calculated table column and the BALANCE
#    formulas in E9/I9 recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> new PERIOD (column A) date serial value.
$dateUpdates = @{
    170 = 44957
    171 = 44985
    172 = 45016
    174 = 45046
    175 = 45077
    176 = 45107
    177 = 45138
    178 = 45169
    179 = 45199
    180 = 45230
    181 = 45260
    182 = 45291
    183 = 45322
    184 = 45351
    185 = 45382
    186 = 45412
    187 = 45443
    188 = 45473
    189 = 45504
    190 = 45535
    191 = 45565
    192 = 45596
    193 = 45626
    194 = 45657
    195 = 45688
    196 = 45716
    197 = 45747
    198 = 45777
    199 = 45808
    200 = 45838
    201 = 45869
    202 = 45900
    203 = 45930
    204 = 45961
    205 = 45991
    206 = 46022
    207 = 46053
    208 = 46081
    209 = 46112
    210 = 46142
    211 = 46173
    212 = 46203
}

foreach ($row in $dateUpdates.Keys) {
    $ws.Cells.Item($row, 1).Value = $dateUpdates[$row]
}

# Newly-posted EARNED credits (column C) for rows 175-179.
$earnedRows = @(175, 176, 177, 178, 179)
foreach ($row in $earnedRows) {
    $ws.Cells.Item($row, 3).Value = 1.25
}

$wb.Save()
